$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24 : new journal entry (09.03.2021) ---
$ws.Cells.Item(24, 2).Value = [DateTime]"2021-03-09"              # B24 Date
$ws.Cells.Item(24, 3).Value = 0.79166666666666663                 # C24 Heure début (19:00)
$ws.Cells.Item(24, 4).Value = 0.82291666666666663                 # D24 Heure fin   (19:45)
$ws.Cells.Item(24, 6).Value = "Ma-20"                              # F24 Module
$ws.Cells.Item(24, 7).Value = "Code"                                # G24 Type
$ws.Cells.Item(24, 8).Value = "Grille"                              # H24 Tâche
$ws.Cells.Item(24, 9).Value = "Maison"                              # I24 Lieu
$ws.Cells.Item(24, 10).Value = "J'ai tester quellque manier de faire une grille et j'ai décider comment j'aillais la faire "  # J24 Descriptif
$ws.Cells.Item(24, 11).Value = "non"                                # K24 Terminer
$ws.Rows.Item(24).RowHeight = 43.2

# --- Row 25 : new journal entry (10.03.2021) ---
$ws.Cells.Item(25, 2).Value = [DateTime]"2021-03-10"              # B25 Date
$ws.Cells.Item(25, 3).Value = 0.5625                                # C25 Heure début (13:30)
$ws.Cells.Item(25, 4).Value = 0.625                                 # D25 Heure fin   (15:00)
$ws.Cells.Item(25, 6).Value = "Ma-20"                              # F25 Module
$ws.Cells.Item(25, 7).Value = "Code"                                # G25 Type
$ws.Cells.Item(25, 8).Value = "Grille"                              # H25 Tâche
$ws.Cells.Item(25, 9).Value = "CPNV"                                # I25 Lieu
$ws.Cells.Item(25, 10).Value = "j'ai finit ma grille et je l'ai afficher dans mon jeux"  # J25 Descriptif
$ws.Cells.Item(25, 11).Value = "Oui"                                # K25 Terminer
$ws.Rows.Item(25).RowHeight = 28.8

# --- Update the sheet view: scroll position and current selection ---
$ws.Activate()
$ws.Range("L25").Select()
